$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets 1-4 ("Potencia Acumulada - SIN (MW)", "Geracao Periodo Medio (MWMed)",
# "Atendimento a Ponta(MW)", "Potencia Incremental - SIN(MW)") all share the
# same layout: add a header label in A1, fix accented labels in A3/A4/A6/A8/A11
# and strip the bold/border header style from A2:A12 (now plain row labels).
# ---------------------------------------------------------------------------
$sheetNames = @(
    "Potencia Acumulada - SIN (MW)",
    "Geracao Periodo Medio (MWMed)",
    "Atendimento a Ponta(MW)",
    "Potencia Incremental - SIN(MW)"
)

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # New header cell A1, copying the existing bold/bordered header format
    # from B1 so it ends up with the identical cell style.
    $ws.Range("B1").Copy() | Out-Null
    $ws.Range("A1").PasteSpecial(-4122) | Out-Null
    $ws.Range("A1").Value = "Fonte/Tecnologia"

    # Fix accented technology labels.
    $ws.Range("A3").Value = "Gás Natural"
    $ws.Range("A4").Value = "Carvão"
    $ws.Range("A6").Value = "Óleos Comb"
    $ws.Range("A8").Value = "Eólica"
    $ws.Range("A11").Value = "Pot. Compl."

    # A2, A5, A7, A9, A10, A12 keep their text but lose the header styling.
    $ws.Range("A2:A12").ClearFormats()
}

# ---------------------------------------------------------------------------
# Sheet 5 ("Emissoes Totais (MtCO2eq)"): add header label, fix accented row
# labels, drop header styling from data rows, and remove the "Teto" row.
# ---------------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item("Emissoes Totais (MtCO2eq)")

$ws5.Range("B1").Copy() | Out-Null
$ws5.Range("A1").PasteSpecial(-4122) | Out-Null
$ws5.Range("A1").Value = "Período"

$ws5.Range("A2").Value = "P.Médio"
$ws5.Range("A3").Value = "P.Crítico"
$ws5.Range("A2:A3").ClearFormats()

# Remove row 4 ("Teto") entirely.
$ws5.Rows.Item(4).Delete() | Out-Null

# ---------------------------------------------------------------------------
# Sheet 6 ("Custo Total (bilhões de R$)"): add header label, rename B1 from
# "Custo" to the text "2015", fix accented labels, strip header styling from
# the data rows and update the cost values.
# ---------------------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("Custo Total (bilhões de R$)")

$ws6.Range("B1").Copy() | Out-Null
$ws6.Range("A1").PasteSpecial(-4122) | Out-Null
$ws6.Range("A1").Value = "Tipo Expansão"

# Keep B1 a text cell (not a number) while reusing its existing style -
# a leading apostrophe forces text entry, matching "2015" used elsewhere.
$ws6.Range("B1").Value = "'2015"

$ws6.Range("A2").Value = "Expansão Centralizada"
$ws6.Range("A3").Value = "Expansão por GD"
$ws6.Range("A2:A3").ClearFormats()

$ws6.Range("B2").Value = 588
$ws6.Range("B3").Value = 99
